# Add bid and comments
# Colors (wdColor / RGB long values, BGR-packed):
#   green  = 00B050 -> RGB(0x00,0xB0,0x50) = 5287936
#   orange = FFC000 -> RGB(0xFF,0xC0,0x00) = 49407
$green  = 5287936
$orange = 49407

$d = $word.ActiveDocument

# --- Hunk 1: Navigation bar paragraph -------------------------------------
# "- Link to all available task page" -> "- " stays plain, rest turns green.
$r = $d.Content
$r.Find.Execute("Link to all available task page") | Out-Null
$r.Font.Color = $green

# "- possibly have search bar to do a basic word matching search? " ->
# "- " stays plain, rest (including trailing space) turns orange.
$r = $d.Content
$r.Find.Execute("possibly have search bar to do a basic word matching search? ") | Out-Null
$r.Font.Color = $orange

# --- Hunk 2: "Individual category pages" paragraph ------------------------
# "- Lists all tasks available for the respective category" turns fully green
# (including the leading "- ").
$r = $d.Content
$r.Find.Execute("- Lists all tasks available for the respective category") | Out-Null
$r.Font.Color = $green

# --- Hunk 3: "Task bidding page" paragraph ---------------------------------
# Everything from "- Or make this reachable..." through the end of the
# paragraph ("...basic task details and minimum bid") turns green, except
# that the line-break just before "- display..." must remain un-colored
# (it stays in its own run). The break character itself occupies the
# position immediately before the found "- display..." text, so we color
# up to (but excluding) that position, then color from the found text's
# start through the paragraph end separately.
$r = $d.Content
$r.Find.Execute("- Or make this reachable from clicking on a task") | Out-Null
$start1 = $r.Start

$r = $d.Content
$r.Find.Execute("- display basic task details and minimum bid") | Out-Null
$start2 = $r.Start
$paraEnd = $r.End

$part1 = $d.Range($start1, $start2 - 1)
$part1.Font.Color = $green

$part2 = $d.Range($start2, $paraEnd)
$part2.Font.Color = $green
